# Make the "Occupation" line conditional on the user being full time
# employed, and collapse the two-run "Total ..." labels in the summary
# tables into single runs (matches the author's template formatting fix).

$d = $word.ActiveDocument

# --- 1. Wrap the Occupation paragraph in a {%p if %} / {%p endif %} block ---

# Locate the paragraph that currently reads "Occupation: {{ user.occup1 }}"
$occupIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Occupation:*") {
        $occupIndex = $i
        break
    }
}

if ($occupIndex -gt 0) {
    $occupParagraph = $d.Paragraphs.Item($occupIndex)

    # Insert a new (empty) paragraph immediately before it, then fill it in.
    $occupParagraph.Range.InsertParagraphBefore()
    $ifParagraph = $d.Paragraphs.Item($occupIndex)
    $ifParagraph.Range.Text = "{%p if user.work == “Full time”%}"

    # The Occupation paragraph has shifted down by one; insert the matching
    # {%p endif %} paragraph right after it.
    $occupParagraph = $d.Paragraphs.Item($occupIndex + 1)
    $occupParagraph.Range.InsertParagraphAfter()
    $endifParagraph = $d.Paragraphs.Item($occupIndex + 2)
    $endifParagraph.Range.Text = "{%p endif %}"
}

# --- 2. Merge the split "Total Value of Assets:" / "{{ total_assets }}" runs ---

$d.Content.Find.Execute(
    "Total Value of Assets: {{ total_assets }}", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "Total Value of Assets: {{ total_assets }}", 2) | Out-Null

# --- 3. Merge the split "Total Amount of Annual Expenses:" / "{{ total_expenses }}" runs ---

$d.Content.Find.Execute(
    "Total Amount of Annual Expenses: {{ total_expenses }}", $false, $false,
    $false, $false, $false, $true, 1, $false,
    "Total Amount of Annual Expenses: {{ total_expenses }}", 2) | Out-Null
